$wb = $excel.ActiveWorkbook

# --- ALC row 128 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H128").Value = 79990
$ws.Range("J128").Value = 79990
$ws.Range("L128").Value = 79990
$ws.Range("N128").Value = -89950

# --- ALC row 129 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1053.4259
$ws.Range("I129").Value = 441.57144
$ws.Range("J129").Value = 1144.5532
$ws.Range("K129").Value = 1324.71432
$ws.Range("L129").Value = 3433.6596
$ws.Range("M129").Value = 3675.28568
$ws.Range("N129").Value = -13433.6596

# --- ALC row 137 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 37038720
$ws.Range("I137").Value = 52632624
$ws.Range("J137").Value = 3196.75
$ws.Range("K137").Value = 157897872
$ws.Range("L137").Value = 9590.25
$ws.Range("M137").Value = -157895322
$ws.Range("N137").Value = -14690.25

# --- ARM row 32 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3576.0154
$ws.Range("I32").Value = 2101.3147
$ws.Range("K32").Value = 2101.3147
$ws.Range("M32").Value = -1814.3147

# --- ARM row 45 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1437.5
$ws.Range("I45").Value = 1300
$ws.Range("J45").Value = 1520
$ws.Range("K45").Value = 1300
$ws.Range("L45").Value = 1520
$ws.Range("M45").Value = -923
$ws.Range("N45").Value = -2274

# --- ARM row 74 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4390.475
$ws.Range("I74").Value = 1389.5
$ws.Range("J74").Value = 11392.75
$ws.Range("K74").Value = 1389.5
$ws.Range("L74").Value = 11392.75
$ws.Range("M74").Value = -515.5
$ws.Range("N74").Value = -13140.75

# --- ARM row 77 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 4390.475
$ws.Range("I77").Value = 1389.5
$ws.Range("J77").Value = 11392.75
$ws.Range("K77").Value = 6947.5
$ws.Range("L77").Value = 56963.75
$ws.Range("M77").Value = -2579.5
$ws.Range("N77").Value = -65699.75

# --- BSM row 107 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 942.4583
$ws.Range("I107").Value = 962.73334
$ws.Range("J107").Value = 908.6667
$ws.Range("K107").Value = 962.73334
$ws.Range("L107").Value = 908.6667
$ws.Range("M107").Value = 957.26666
$ws.Range("N107").Value = -4748.6667

# --- BSM row 134 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3133.0645
$ws.Range("I134").Value = 2223
$ws.Range("J134").Value = 5357.6665
$ws.Range("K134").Value = 6669
$ws.Range("L134").Value = 16072.9995
$ws.Range("M134").Value = -4134
$ws.Range("N134").Value = -21142.9995

# --- CRP row 99 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 15626500
$ws.Range("I99").Value = 15626500
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 15626500
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -15625002
$ws.Range("N99").ClearContents()

# --- CRP row 126 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 15626500
$ws.Range("I126").Value = 15626500
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 46879500
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -46877030
$ws.Range("N126").ClearContents()

# --- CUL row 34 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 3479.8
$ws.Range("J34").Value = 4199.75
$ws.Range("L34").Value = 12599.25
$ws.Range("N34").Value = -12767.25

# --- CUL row 39 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 7655.5
$ws.Range("J39").Value = 8109.0303
$ws.Range("L39").Value = 24327.0909
$ws.Range("N39").Value = -24915.0909

# --- CUL row 51 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 8000
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 8000
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 24000
$ws.Range("N51").Value = -24920
$ws.Range("M51").ClearContents()

# --- CUL row 55 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 36700
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 36700
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 110100
$ws.Range("N55").Value = -110454
$ws.Range("M55").ClearContents()

# --- CUL row 57 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 6000
$ws.Range("J57").Value = 6000
$ws.Range("L57").Value = 18000
$ws.Range("N57").Value = -19118

# --- CUL row 58 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 2082.2
$ws.Range("J58").Value = 1268.6666
$ws.Range("L58").Value = 3805.9998
$ws.Range("N58").Value = -4061.9998

# --- CUL row 64 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 3400
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 3400
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 10200
$ws.Range("N64").Value = -10740
$ws.Range("M64").ClearContents()

# --- CUL row 67 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H67").Value = 3400
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 3400
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 10200
$ws.Range("N67").Value = -12072
$ws.Range("M67").ClearContents()

# --- CUL row 70 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 1914
$ws.Range("I70").Value = 1262.4
$ws.Range("J70").Value = 3000
$ws.Range("K70").Value = 3787.2
$ws.Range("L70").Value = 9000
$ws.Range("M70").Value = -3472.2
$ws.Range("N70").Value = -9630

# --- CUL row 73 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H73").Value = 1914
$ws.Range("I73").Value = 1262.4
$ws.Range("J73").Value = 3000
$ws.Range("K73").Value = 3787.2
$ws.Range("L73").Value = 9000
$ws.Range("M73").Value = -2695.2
$ws.Range("N73").Value = -11184

# --- CUL row 76 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 2999
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

# --- CUL row 79 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H79").Value = 2999
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

# --- GSM row 102 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2011.7693
$ws.Range("I102").Value = 1809.8
$ws.Range("J102").Value = 2138
$ws.Range("K102").Value = 1809.8
$ws.Range("L102").Value = 2138
$ws.Range("M102").Value = -187.8
$ws.Range("N102").Value = -5382

# --- GSM row 132 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2651.309
$ws.Range("I132").Value = 2112.5293
$ws.Range("J132").Value = 3523.6191
$ws.Range("K132").Value = 6337.5879
$ws.Range("L132").Value = 10570.8573
$ws.Range("M132").Value = -3807.5879
$ws.Range("N132").Value = -15630.8573

# --- LTW row 136 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3331.6047
$ws.Range("I136").Value = 1848.3103
$ws.Range("J136").Value = 6404.143
$ws.Range("K136").Value = 5544.9309
$ws.Range("L136").Value = 19212.429
$ws.Range("M136").Value = -2994.9309
$ws.Range("N136").Value = -24312.429

# --- WVR row 132 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 38465316
$ws.Range("I132").Value = 50003790
$ws.Range("J132").Value = 3732.6667
$ws.Range("K132").Value = 150011370
$ws.Range("L132").Value = 11198.0001
$ws.Range("M132").Value = -150008840
$ws.Range("N132").Value = -16258.0001
